$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados..." timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 00:52"

# --- Country name reorder (Burkina Faso's case count overtook Letonia's,
#     so it moved up in the ranking, ahead of Letonia and Birmania) ---
$ws.Cells.Item(156, 1).Value = "Burkina Faso"
$ws.Cells.Item(157, 1).Value = "Letonia"
$ws.Cells.Item(158, 1).Value = "Birmania"

# --- Updated per-country statistics ---
# columns B..H = Casos totales, Nuevos casos, Casos activos, Recuperados,
#                Casos criticos, Muertes hoy, Muertes
# each entry: row number, then the 7 values for B:H

$rowData = @(
    @(4,   6457446, 28306, 3720625, 2543615, 0, 388, 193206),  # Estados Unidos
    @(6,   4137521, 14521, 3317227,  693644, 0, 420, 126650),  # Brasil
    @(9,    666521,  8065,  518229,  126880, 0, 256,  21412),  # Colombia
    @(13,   478792,  6986,  349132,  119801, 0, 120,   9859),  # Argentina
    @(34,    99863,   151,   78108,   16225, 0,  19,   5530),  # Egipto
    @(43,    77683,   202,   66131,    8700, 0,   7,   2852),  # Guatemala
    @(48,    71419,   543,   62076,    7986, 0,   8,   1357),  # Japon
    @(55,    55005,   100,   43013,   10935, 0,   3,   1057),  # Nigeria
    @(56,    53289,  1124,   42006,   10855, 0,   8,    428),  # Venezuela
    @(74,    26279,    72,   22467,    3059, 0,   5,    753),  # Australia
    @(84,    17089,    39,   12157,    4256, 0,   5,    676),  # Bulgaria
    @(88,    13437,    30,    6730,    5874, 0,   1,    833),  # Sudan
    @(92,    11388,    92,    9348,    1776, 0,   0,    264),  # Noruega
    @(95,     9798,   149,    8928,     809, 0,   0,     61),  # Guinea
    @(139,    2449,    21,    1206,    1229, 0,   0,     14),  # Aruba
    @(151,    1679,    10,    1459,     175, 0,   0,     45),  # Uruguay
    @(155,    1468,     9,     954,     468, 0,   0,     46),  # Guyana
    @(156,    1452,     5,    1103,     294, 0,   0,     55),  # Burkina Faso (new)
    @(157,    1428,     3,    1187,     206, 0,   0,     35),  # Letonia (shifted)
    @(158,    1419,   166,     385,    1026, 0,   1,      8),  # Birmania (shifted)
    @(166,    1039,     5,     918,      42, 0,   2,     79)   # Republica del Chad
)

foreach ($row in $rowData) {
    $r = $row[0]
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $row[$c + 1]
    }
}
